$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously contained age-like placeholder values (1-6) in the
# "Revenue" column by mistake. Replace them with the correct revenue figures.
$ws.Range("A2").Value = 50000
$ws.Range("A3").Value = 13456
$ws.Range("A4").Value = 134567
$ws.Range("A5").Value = 76544
$ws.Range("A6").Value = 357373
$ws.Range("A7").Value = 263573

# Move/update the active selection to F3, matching the author's last
# on-screen selection before saving.
$ws.Range("F3").Select()
